$d = $word.ActiveDocument

# The "Langages" line currently sits after "MLOps" (and right before
# "Bases de donnees"); it needs to move up so it becomes the first
# line of the COMPETENCES TECHNIQUES block, i.e. right before the
# "Visualisation : tableau" paragraph.

$langagesText = "Langages : python, matlab, c, c++"
$anchorText = "Visualisation : tableau"

$langagesIndex = -1
$anchorIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r", "`a")
    if ($langagesIndex -eq -1 -and $t -eq $langagesText) {
        $langagesIndex = $i
    }
    if ($anchorIndex -eq -1 -and $t -eq $anchorText) {
        $anchorIndex = $i
    }
}

if ($langagesIndex -ne -1 -and $anchorIndex -ne -1 -and $langagesIndex -ne $anchorIndex) {
    # Paragraph objects/indices shift as soon as the document is
    # mutated, so do the removal/insertion in an order that never
    # needs a stale index: handle whichever paragraph sits later in
    # the document first, since acting on it cannot move the index of
    # anything that comes before it.
    if ($langagesIndex -gt $anchorIndex) {
        # Remove the original "Langages" paragraph (including its
        # paragraph mark) first - it is after the anchor, so the
        # anchor's index is unaffected.
        $d.Paragraphs.Item($langagesIndex).Range.Delete()
        # Re-insert the text right before the anchor paragraph; the
        # new paragraph inherits the anchor's paragraph formatting.
        $d.Paragraphs.Item($anchorIndex).Range.InsertBefore($langagesText + "`r")
    } else {
        # Anchor is after "Langages": insert first (this shifts the
        # anchor index by one, but not the still-untouched original
        # "Langages" paragraph, which stays at $langagesIndex).
        $d.Paragraphs.Item($anchorIndex).Range.InsertBefore($langagesText + "`r")
        $d.Paragraphs.Item($langagesIndex).Range.Delete()
    }
}
